$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update H4 to "Done" and I4 to "In-Progress" (label progress tracking update)
$ws.Range("H4").Value = "Done"
$ws.Range("I4").Value = "In-Progress"

# Update the active selection to H8
$ws.Range("H8").Select()
